$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Tornado Damages Connecticut River Walk in Springfield, Mass.'
$ws.Range('B2').Value = '2011-06-08T10:24:00UTC'
$ws.Range('C2').Value = 7
$ws.Range('D2').Value = 'day_2_to_30'
$ws.Range('E2').Value = 'https://web.archive.org/web/20110611080404/http://community.railstotrails.org/blogs/trailblog/archive/2011/06/03/tornado-damages-springfield-riverwalk-in-massachusetts.aspx'

$ws.Range('A3').Value = 'Gov. Deval Patrick requests federal disaster aid for tornadoes'
$ws.Range('B3').Value = '2011-06-11T00:00:00UTC'
$ws.Range('C3').Value = 10
$ws.Range('D3').Value = 'day_2_to_30'
$ws.Range('E3').Value = 'http://www.bostonherald.com/news/regional/view/2011_0611gov_patrick_requests_federal_disaster_aid_for_tornadoes/'

$ws.Range('A4').Value = 'Massachusetts Tornadoes: At Least 3 Dead; Springfield, Westfield Hard Hit'
$ws.Range('B4').Value = '2011-06-02T00:00:00UTC'
$ws.Range('C4').Value = 1
$ws.Range('D4').Value = 'day_1'
$ws.Range('E4').Value = 'https://abcnews.go.com/US/massachusetts-rocked-multiple-tornadoes-including-city-springfield-possibly/story?id=13737522&sms_ss=facebook&at_xt=4de6f2afa93a9b3a,0'

$ws.Range('A5').Value = 'Victims lament Massachusetts tornadoes’ huge toll'
$ws.Range('B5').Value = '2011-06-03T00:00:00UTC'
$ws.Range('C5').Value = 2
$ws.Range('D5').Value = 'day_2_to_30'
$ws.Range('E5').Value = 'http://www.boston.com/news/local/massachusetts/articles/2011/06/03/victims_lament_massachusetts_tornadoes_huge_toll/'

$ws.Range('A6').Value = '65 years ago, nightmare tornado killed 94 Central Mass.'
$ws.Range('B6').Value = '2018-06-09T07:00:00UTC'
$ws.Range('C6').Value = 2565
$ws.Range('D6').Value = 'day_31_beyond'
$ws.Range('E6').Value = 'https://www.telegram.com/news/20180609/65-years-ago-nightmare-tornado-killed-94-central-mass'

$ws.Range('A7').Value = 'Springfield tornado death toll at 3; Massachusetts Congressional delegation appeals to President Obama for disaster declaration'
$ws.Range('B7').Value = '2011-06-03T06:00:07UTC'
$ws.Range('C7').Value = 2
$ws.Range('D7').Value = 'day_2_to_30'
$ws.Range('E7').Value = 'http://www.masslive.com/news/index.ssf/2011/06/springfield_tornado_death_toll.html'

$ws.Range('A8').Value = 'Massachusetts digs out after three tornadoes kill 3, hurt 200'
$ws.Range('B8').Value = '2011-06-01T00:00:00UTC'
$ws.Range('C8').Value = 0
$ws.Range('D8').Value = 'day_0'
$ws.Range('E8').Value = 'http://www.myfoxboston.com/dpp/weather/tornado-damage-reported-in-springfield-25-apx-20110601'

$ws.Range('A9').Value = '40 years ago today, a tornado ripped through West Stockbridge'
$ws.Range('B9').Value = '2013-08-28T06:05:12UTC'
$ws.Range('C9').Value = 819
$ws.Range('D9').Value = 'day_31_beyond'
$ws.Range('E9').Value = 'https://www.berkshireeagle.com/stories/40-years-ago-today-a-tornado-ripped-through-west-stockbridge,400213'

$ws.Range('A10').Value = 'Insurance claims for tornadoes hit $140m'
$ws.Range('B10').Value = '2011-06-18T00:00:00UTC'
$ws.Range('C10').Value = 17
$ws.Range('D10').Value = 'day_2_to_30'
$ws.Range('E10').Value = 'https://web.archive.org/web/20120406092705/http://articles.boston.com/2011-06-18/news/29674671_1_recovery-centers-insurance-claims-fema'

$ws.Range('A11').Value = 'Brother’s tornado death rocks West Springfield family'
$ws.Range('B11').Value = '2011-06-02T00:00:00UTC'
$ws.Range('C11').Value = 1
$ws.Range('D11').Value = 'day_1'
$ws.Range('E11').Value = 'http://news.bostonherald.com/news/regional/view/2011_0602brothers_tornado_death_rocks_west_springfield_family/'

$ws.Range('A12').Value = 'Mother gives up life to save her daughter during tornado'
$ws.Range('B12').Value = '2011-06-03T00:00:00UTC'
$ws.Range('C12').Value = 2
$ws.Range('D12').Value = 'day_2_to_30'
$ws.Range('E12').Value = 'https://web.archive.org/web/20110603111417/http://www.boston.com/news/weather/articles/2011/06/03/mother_gives_up_life_to_save_her_daughter_during_tornado/'

$ws.Range('A13').Value = 'Cathedral High School: A review of the latest stories on future of Springfield school'
$ws.Range('B13').Value = '2015-01-31T12:05:06UTC'
$ws.Range('C13').Value = 1340
$ws.Range('D13').Value = 'day_31_beyond'
$ws.Range('E13').Value = 'http://www.masslive.com/news/index.ssf/2015/01/several_days_of_news_on_the_future_of_cathedral_follow_the_latest.html'

$ws.Range('A14').Value = 'Red Cross Helping After Tornadoes Tear Through Massachusetts'
$ws.Range('B14').Value = '2011-06-02T13:00:00UTC'
$ws.Range('C14').Value = 1
$ws.Range('D14').Value = 'day_1'
$ws.Range('E14').Value = 'http://www.prnewswire.com/news-releases/red-cross-helping-after-tornadoes-tear-through-massachusetts-123031908.html'

$ws.Range('A15').Value = 'Springfield Massachusetts tornado and other June 1st storms on the Storm Time Machine'
$ws.Range('B15').Value = '2011-06-01T00:00:00UTC'
$ws.Range('C15').Value = 0
$ws.Range('D15').Value = 'day_0'
$ws.Range('E15').Value = 'http://www.stormtimemachine.com/20110601.html'

$ws.Range('A16').Value = 'Hundreds still homeless after tornadoes'
$ws.Range('B16').Value = '1-01-01T00:00:00UTC'
$ws.Range('C16').Value = 'unknown'
$ws.Range('D16').Value = 'unknown'
$ws.Range('E16').Value = 'https://web.archive.org/web/20120322001638/http://www.wwlp.com/dpp/news/local/hampden/hundreds-still-homeless-after-tornadoes'

$ws.Range('A17').Value = 'National Weather Service Text Product Display'
$ws.Range('B17').Value = '1-01-01T00:00:00UTC'
$ws.Range('C17').Value = 'unknown'
$ws.Range('D17').Value = 'unknown'
$ws.Range('E17').Value = 'http://forecast.weather.gov/product.php?site%3DNWS%26issuedby%3DGYX%26product%3DPNS%26format%3DCI%26version%3D2%26glossary%3D0'

$ws.Range('A18').Value = 'In Brimfield, terror -- and death -- arrived in a cloud laden with debris'
$ws.Range('B18').Value = '1-01-01T00:00:00UTC'
$ws.Range('C18').Value = 'unknown'
$ws.Range('D18').Value = 'unknown'
$ws.Range('E18').Value = 'https://web.archive.org/web/20110609020449/http://www.boston.com/news/local/breaking_news/2011/06/in_brimfield_te.html'

$ws.Range('A19').Value = 'Storm Prediction Center Severe Thunderstorm Watch 411'
$ws.Range('B19').Value = '1-01-01T00:00:00UTC'
$ws.Range('C19').Value = 'unknown'
$ws.Range('D19').Value = 'unknown'
$ws.Range('E19').Value = 'http://www.spc.noaa.gov/products/watch/ww0411.html'

$ws.Range('A20').Value = 'Tornado History Project: 20110601.25.1'
$ws.Range('B20').Value = '1-01-01T00:00:00UTC'
$ws.Range('C20').Value = 'unknown'
$ws.Range('D20').Value = 'unknown'
$ws.Range('E20').Value = 'http://www.tornadohistoryproject.com/tornado/20110601.25.1'

$ws.Range('A21').Value = 'Tornado History Project: Massachusetts'
$ws.Range('B21').Value = '1-01-01T00:00:00UTC'
$ws.Range('C21').Value = 'unknown'
$ws.Range('D21').Value = 'unknown'
$ws.Range('E21').Value = 'http://www.tornadohistoryproject.com/tornado/Massachusetts/map'

$ws.Range('A22').Value = 'News and Information from Northampton, MA by the Daily Hampshire Gazette'
$ws.Range('B22').Value = '1-01-01T00:00:00UTC'
$ws.Range('C22').Value = 'unknown'
$ws.Range('D22').Value = 'unknown'
$ws.Range('E22').Value = 'http://www.gazettenet.com/2011/06/02/governor-patrick-sen-kerry-survey-destruction-in-western-massachusetts'

$ws.Range('A23').Value = 'Storm Prediction Center Severe Thunderstorm Watch 410'
$ws.Range('B23').Value = '1-01-01T00:00:00UTC'
$ws.Range('C23').Value = 'unknown'
$ws.Range('D23').Value = 'unknown'
$ws.Range('E23').Value = 'http://www.spc.noaa.gov/products/watch/ww0410.html'

$ws.Range('A24').Value = 'Massachusetts Severe Storms and Tornadoes'
$ws.Range('B24').Value = '1-01-01T00:00:00UTC'
$ws.Range('C24').Value = 'unknown'
$ws.Range('D24').Value = 'unknown'
$ws.Range('E24').Value = 'http://www.fema.gov/news/eventcounties.fema?id=14733'

$ws.Range('A25').Value = 'Annual Fatal Tornado Summaries'
$ws.Range('B25').Value = '1-01-01T00:00:00UTC'
$ws.Range('C25').Value = 'unknown'
$ws.Range('D25').Value = 'unknown'
$ws.Range('E25').Value = 'https://web.archive.org/web/20110613140257/http://www.spc.noaa.gov/climo/torn/fataltorn.html'

$ws.Range('A26').Value = 'Helping Homeless Tornado Victims'
$ws.Range('B26').Value = '1-01-01T00:00:00UTC'
$ws.Range('C26').Value = 'unknown'
$ws.Range('D26').Value = 'unknown'
$ws.Range('E26').Value = 'https://web.archive.org/web/20111002155822/http://www.wggb.com/story/14943481/helping-homeless-tornado-victims'
